# winch functions: cyrpus digial 1-2, cyrpus analog 1
#
# - "winch lock" (Outputs!E14) is re-cased to "Winch Lock"
# - "Other Inputs" rows 11-17 gain the new Cyprus winch-lock entries:
#     row 11: CYPRUS STUFF / Digital / D / 1 / Hold to fire winch
#     row 12:                         D / 2 / Push to lock / unlock
#     row 13-16:                      D / 3..6
#     row 17:        Dial   / A / 1 / Change winch speed
# - selection/active-cell bookmarks on the Outputs and Other Inputs sheets move

$wb = $excel.ActiveWorkbook

$wsOut = $wb.Worksheets.Item("Outputs")
$wsOut.Range("E14").Value = "Winch Lock"

$wsOther = $wb.Worksheets.Item("Other Inputs")

# Row 11 - fill variable name, class, channel and function first...
$wsOther.Range("A11").Value = "CYPRUS STUFF"
$wsOther.Range("B11").Value = "Digital"
$wsOther.Range("D11").Value = 1
$wsOther.Range("E11").Value = "Hold to fire winch"

# Row 12 - channel + function
$wsOther.Range("D12").Value = 2
$wsOther.Range("E12").Value = "Push to lock / unlock"

# ...then go back and fill the "Slot" column (D) for all six digital channels
$wsOther.Range("C11").Value = "D"
$wsOther.Range("C12").Value = "D"
$wsOther.Range("C13").Value = "D"
$wsOther.Range("C14").Value = "D"
$wsOther.Range("C15").Value = "D"
$wsOther.Range("C16").Value = "D"

# remaining channel numbers
$wsOther.Range("D13").Value = 3
$wsOther.Range("D14").Value = 4
$wsOther.Range("D15").Value = 5
$wsOther.Range("D16").Value = 6

# Row 17 - the analog dial entry
$wsOther.Range("C17").Value = "A"
$wsOther.Range("B17").Value = "Dial"
$wsOther.Range("D17").Value = 1
$wsOther.Range("E17").Value = "Change winch speed"

# Restore the recorded selections (view state) on each sheet
$wsOut.Activate()
$wsOut.Range("E15").Select()

$wsOther.Activate()
$wsOther.Range("B18").Select()
